$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at the top, shifting existing data down
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "Header row description"
$ws.Range("B1").Value = "Value"

[void]$ws.Range("A20").Select()
